$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1957.2428
$ws.Range("J17").Value = 1957.2428
$ws.Range("L17").Value = 5871.7284
$ws.Range("N17").Value = -6207.7284
$ws.Range("H125").Value = 1888.6
$ws.Range("I125").Value = 1527.8572
$ws.Range("J125").Value = 2204.25
$ws.Range("K125").Value = 13750.7148
$ws.Range("L125").Value = 19838.25
$ws.Range("M125").Value = -11290.7148
$ws.Range("N125").Value = -24758.25
$ws.Range("H137").Value = 2573350
$ws.Range("I137").Value = 4533752.5
$ws.Range("J137").Value = 9746.691999999999
$ws.Range("K137").Value = 13601257.5
$ws.Range("L137").Value = 29240.076
$ws.Range("M137").Value = -13598707.5
$ws.Range("N137").Value = -34340.076
$ws.Range("H141").Value = 3553.5
$ws.Range("I141").Value = 1932.3636
$ws.Range("J141").Value = 5534.8887
$ws.Range("K141").Value = 5797.0908
$ws.Range("L141").Value = 16604.6661
$ws.Range("M141").Value = -617.0907999999999
$ws.Range("N141").Value = -26964.6661

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9259.986999999999
$ws.Range("I32").Value = 8171.4604
$ws.Range("J32").Value = 25805.6
$ws.Range("K32").Value = 8171.4604
$ws.Range("L32").Value = 25805.6
$ws.Range("M32").Value = -7884.4604
$ws.Range("N32").Value = -26379.6
$ws.Range("H63").Value = 3320.0625
$ws.Range("I63").Value = 2315.5
$ws.Range("J63").Value = 4994.3335
$ws.Range("K63").Value = 2315.5
$ws.Range("L63").Value = 4994.3335
$ws.Range("M63").Value = -1629.5
$ws.Range("N63").Value = -6366.3335
$ws.Range("H66").Value = 3320.0625
$ws.Range("I66").Value = 2315.5
$ws.Range("J66").Value = 4994.3335
$ws.Range("K66").Value = 11577.5
$ws.Range("L66").Value = 24971.6675
$ws.Range("M66").Value = -8145.5
$ws.Range("N66").Value = -31835.6675
$ws.Range("H122").Value = 1810.85
$ws.Range("I122").Value = 1531
$ws.Range("J122").Value = 2330.5715
$ws.Range("K122").Value = 4593
$ws.Range("L122").Value = 6991.7145
$ws.Range("M122").Value = -2143
$ws.Range("N122").Value = -11891.7145
$ws.Range("H132").Value = 14707885
$ws.Range("I132").Value = 20834846
$ws.Range("J132").Value = 3178.4
$ws.Range("K132").Value = 62504538
$ws.Range("L132").Value = 9535.200000000001
$ws.Range("M132").Value = -62502008
$ws.Range("N132").Value = -14595.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 28370
$ws.Range("J35").Value = 28370
$ws.Range("L35").Value = 28370
$ws.Range("N35").Value = -28990
$ws.Range("H82").Value = 32051.4
$ws.Range("I82").Value = 5128.5
$ws.Range("K82").Value = 5128.5
$ws.Range("M82").Value = -4745.5
$ws.Range("H85").Value = 32051.4
$ws.Range("I85").Value = 5128.5
$ws.Range("K85").Value = 5128.5
$ws.Range("M85").Value = -3802.5
$ws.Range("H134").Value = 2493.3416
$ws.Range("I134").Value = 2060
$ws.Range("J134").Value = 4280.875
$ws.Range("K134").Value = 6180
$ws.Range("L134").Value = 12842.625
$ws.Range("M134").Value = -3645
$ws.Range("N134").Value = -17912.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 22719.5
$ws.Range("I41").Value = 4059
$ws.Range("J41").Value = 28939.666
$ws.Range("K41").Value = 4059
$ws.Range("L41").Value = 28939.666
$ws.Range("M41").Value = -3631
$ws.Range("N41").Value = -29795.666
$ws.Range("H50").Value = 32343.8
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 39179.75
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 39179.75
$ws.Range("M50").Value = -4375
$ws.Range("N50").Value = -40429.75
$ws.Range("H51").Value = 41702800
$ws.Range("I51").Value = 500000000
$ws.Range("J51").Value = 39419.91
$ws.Range("K51").Value = 500000000
$ws.Range("L51").Value = 39419.91
$ws.Range("M51").Value = -499999264
$ws.Range("N51").Value = -40891.91
$ws.Range("H58").Value = 2567.2222
$ws.Range("I58").Value = 1624
$ws.Range("J58").Value = 3443.0715
$ws.Range("K58").Value = 1624
$ws.Range("L58").Value = 3443.0715
$ws.Range("M58").Value = -1421
$ws.Range("N58").Value = -3849.0715
$ws.Range("H60").Value = 20100
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 20100
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 20100
$ws.Range("N60").Value = -21122
$ws.Range("H61").Value = 41702800
$ws.Range("I61").Value = 500000000
$ws.Range("J61").Value = 39419.91
$ws.Range("K61").Value = 500000000
$ws.Range("L61").Value = 39419.91
$ws.Range("M61").Value = -499999652
$ws.Range("N61").Value = -40115.91
$ws.Range("H99").Value = 1992.16
$ws.Range("I99").Value = 1945.8667
$ws.Range("J99").Value = 2061.6
$ws.Range("K99").Value = 1945.8667
$ws.Range("L99").Value = 2061.6
$ws.Range("M99").Value = -447.8667
$ws.Range("N99").Value = -5057.6
$ws.Range("H122").Value = 64177
$ws.Range("I122").Value = 64177
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 192531
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -190081
$ws.Range("H126").Value = 1992.16
$ws.Range("I126").Value = 1945.8667
$ws.Range("J126").Value = 2061.6
$ws.Range("K126").Value = 5837.6001
$ws.Range("L126").Value = 6184.799999999999
$ws.Range("M126").Value = -3367.6001
$ws.Range("N126").Value = -11124.8
$ws.Range("H131").Value = 59326
$ws.Range("J131").Value = 59326
$ws.Range("L131").Value = 59326
$ws.Range("H136").Value = 2567.2222
$ws.Range("I136").Value = 1624
$ws.Range("J136").Value = 3443.0715
$ws.Range("K136").Value = 4872
$ws.Range("L136").Value = 10329.2145
$ws.Range("M136").Value = -2322
$ws.Range("N136").Value = -15429.2145
$ws.Range("N131").Value = -69406
$ws.Range("M60").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6156.3335
$ws.Range("I3").Value = 2300
$ws.Range("J3").Value = 8084.5
$ws.Range("K3").Value = 6900
$ws.Range("L3").Value = 24253.5
$ws.Range("M3").Value = -6788
$ws.Range("N3").Value = -24477.5
$ws.Range("H68").Value = 1341.6451
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1341.6451
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4024.9353
$ws.Range("N68").Value = -5646.9353
$ws.Range("H71").Value = 1341.6451
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1341.6451
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 12074.8059
$ws.Range("N71").Value = -20186.8059
$ws.Range("H107").Value = 11698.3
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 11698.3
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 35094.89999999999
$ws.Range("N107").Value = -38934.89999999999
$ws.Range("H110").Value = 1000
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H111").Value = 4000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 12000
$ws.Range("N111").Value = -18134
$ws.Range("H125").Value = 1672060.5
$ws.Range("J125").Value = 6072.6665
$ws.Range("L125").Value = 18217.9995
$ws.Range("N125").Value = -28057.9995
$ws.Range("H130").Value = 61285.6
$ws.Range("J130").Value = 1966
$ws.Range("L130").Value = 5898
$ws.Range("N130").Value = -15938
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("M107").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("M111").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1103
$ws.Range("I102").Value = 964.8
$ws.Range("J102").Value = 1333.3334
$ws.Range("K102").Value = 964.8
$ws.Range("L102").Value = 1333.3334
$ws.Range("M102").Value = 657.2
$ws.Range("N102").Value = -4577.3334
$ws.Range("H122").Value = 1528.5714
$ws.Range("I122").Value = 1533.3334
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4600.0002
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2150.0002
$ws.Range("H126").Value = 9392.267
$ws.Range("I126").Value = 22724.8
$ws.Range("J126").Value = 2726
$ws.Range("K126").Value = 68174.39999999999
$ws.Range("L126").Value = 8178
$ws.Range("M126").Value = -65704.39999999999
$ws.Range("N126").Value = -13118
$ws.Range("H132").Value = 40004770
$ws.Range("I132").Value = 66671496
$ws.Range("J132").Value = 4680.5
$ws.Range("K132").Value = 200014488
$ws.Range("L132").Value = 14041.5
$ws.Range("M132").Value = -200011958
$ws.Range("N132").Value = -19101.5
$ws.Range("N122").Value = -9400

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3015.9
$ws.Range("I132").Value = 2285.6128
$ws.Range("K132").Value = 6856.8384
$ws.Range("M132").Value = -4326.8384
